$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1740
$ws.Range("J40").Value = 1750
$ws.Range("L40").Value = 1750
$ws.Range("N40").Value = -2100

# Row 137
$ws.Range("H137").Value = 15635975
$ws.Range("I137").Value = 3472809
$ws.Range("J137").Value = 52125476
$ws.Range("K137").Value = 10418427
$ws.Range("L137").Value = 156376428
$ws.Range("M137").Value = -10415877
$ws.Range("N137").Value = -156381528

# Row 138
$ws.Range("H138").Value = 2977.24
$ws.Range("I138").Value = 3376.5
$ws.Range("J138").Value = 2851.158
$ws.Range("K138").Value = 10129.5
$ws.Range("L138").Value = 8553.474
$ws.Range("M138").Value = -4989.5
$ws.Range("N138").Value = -18833.474

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 10429.083
$ws.Range("I2").Value = 13681.296
$ws.Range("J2").Value = 672.44446
$ws.Range("K2").Value = 13681.296
$ws.Range("L2").Value = 672.44446
$ws.Range("M2").Value = -13568.296
$ws.Range("N2").Value = -898.44446

# Row 44
$ws.Range("H44").Value = 20499.8
$ws.Range("J44").Value = 20499.8
$ws.Range("L44").Value = 20499.8
$ws.Range("N44").Value = -21475.8

# Row 61
$ws.Range("H61").Value = 2793049.5
$ws.Range("I61").Value = 1603813.8
$ws.Range("J61").Value = 5885063
$ws.Range("K61").Value = 1603813.8
$ws.Range("L61").Value = 5885063
$ws.Range("M61").Value = -1603601.8
$ws.Range("N61").Value = -5885487

# Row 80
$ws.Range("H80").Value = 26748.875
$ws.Range("J80").Value = 26748.875
$ws.Range("L80").Value = 26748.875
$ws.Range("N80").Value = -28744.875

# Row 83
$ws.Range("H83").Value = 26748.875
$ws.Range("J83").Value = 26748.875
$ws.Range("L83").Value = 80246.625
$ws.Range("N83").Value = -90230.625

# Row 88
$ws.Range("H88").Value = 4415.6
$ws.Range("J88").Value = 4706.222
$ws.Range("L88").Value = 4706.222
$ws.Range("N88").Value = -5518.222

# Row 91
$ws.Range("H91").Value = 4415.6
$ws.Range("J91").Value = 4706.222
$ws.Range("L91").Value = 4706.222
$ws.Range("N91").Value = -7514.222

# Row 116
$ws.Range("H116").Value = 10429.083
$ws.Range("I116").Value = 13681.296
$ws.Range("J116").Value = 672.44446
$ws.Range("K116").Value = 13681.296
$ws.Range("L116").Value = 672.44446
$ws.Range("M116").Value = -11387.296
$ws.Range("N116").Value = -5260.44446

# Row 122
$ws.Range("H122").Value = 1705.0416
$ws.Range("I122").Value = 1902.0769
$ws.Range("J122").Value = 1472.1818
$ws.Range("K122").Value = 5706.2307
$ws.Range("L122").Value = 4416.5454
$ws.Range("M122").Value = -3256.2307
$ws.Range("N122").Value = -9316.545399999999

# Row 136
$ws.Range("H136").Value = 2793049.5
$ws.Range("I136").Value = 1603813.8
$ws.Range("J136").Value = 5885063
$ws.Range("K136").Value = 4811441.4
$ws.Range("L136").Value = 17655189
$ws.Range("M136").Value = -4808891.4
$ws.Range("N136").Value = -17660289

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 10429.083
$ws.Range("I3").Value = 13681.296
$ws.Range("J3").Value = 672.44446
$ws.Range("K3").Value = 13681.296
$ws.Range("L3").Value = 672.44446
$ws.Range("M3").Value = -13567.296
$ws.Range("N3").Value = -900.44446

# Row 86
$ws.Range("H86").Value = 1870.909
$ws.Range("I86").Value = 1903.2211
$ws.Range("K86").Value = 1903.2211
$ws.Range("M86").Value = -780.2211

# Row 89
$ws.Range("H89").Value = 1870.909
$ws.Range("I89").Value = 1903.2211
$ws.Range("K89").Value = 9516.1055
$ws.Range("M89").Value = -3900.1055

# Row 105
$ws.Range("H105").Value = 1520
$ws.Range("I105").Value = 1614.2858
$ws.Range("J105").Value = 1300
$ws.Range("K105").Value = 1614.2858
$ws.Range("L105").Value = 1300
$ws.Range("M105").Value = 132.7141999999999
$ws.Range("N105").Value = -4794

# Row 122
$ws.Range("H122").Value = 35392
$ws.Range("J122").Value = 35392
$ws.Range("L122").Value = 35392
$ws.Range("N122").Value = -45192

# Row 134
$ws.Range("H134").Value = 12119757
$ws.Range("I134").Value = 14360647
$ws.Range("J134").Value = 4052555
$ws.Range("K134").Value = 43081941
$ws.Range("L134").Value = 12157665
$ws.Range("M134").Value = -43079406
$ws.Range("N134").Value = -12162735

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 939834.4
$ws.Range("I31").Value = 2875285.8
$ws.Range("J31").Value = 4366.1167
$ws.Range("K31").Value = 2875285.8
$ws.Range("L31").Value = 4366.1167
$ws.Range("M31").Value = -2874990.8
$ws.Range("N31").Value = -4956.1167

# Row 34
$ws.Range("H34").Value = 939834.4
$ws.Range("I34").Value = 2875285.8
$ws.Range("J34").Value = 4366.1167
$ws.Range("K34").Value = 2875285.8
$ws.Range("L34").Value = 4366.1167
$ws.Range("M34").Value = -2875083.8
$ws.Range("N34").Value = -4770.1167

# Row 58
$ws.Range("H58").Value = 1140743.4
$ws.Range("I58").Value = 5729.095
$ws.Range("J58").Value = 2395233
$ws.Range("K58").Value = 5729.095
$ws.Range("L58").Value = 2395233
$ws.Range("M58").Value = -5526.095
$ws.Range("N58").Value = -2395639

# Row 122
$ws.Range("H122").Value = 4213.7646
$ws.Range("I122").Value = 5849.8096
$ws.Range("K122").Value = 17549.4288
$ws.Range("M122").Value = -15099.4288

# Row 136
$ws.Range("H136").Value = 1140743.4
$ws.Range("I136").Value = 5729.095
$ws.Range("J136").Value = 2395233
$ws.Range("K136").Value = 17187.285
$ws.Range("L136").Value = 7185699
$ws.Range("M136").Value = -14637.285
$ws.Range("N136").Value = -7190799

$ws = $wb.Worksheets.Item("CUL")
# Row 127
$ws.Range("H127").Value = 1150
$ws.Range("I127").Value = 450
$ws.Range("J127").Value = 1250
$ws.Range("K127").Value = 1350
$ws.Range("L127").Value = 3750
$ws.Range("M127").Value = 3610
$ws.Range("N127").Value = -13670

# Row 131
$ws.Range("H131").Value = 885.12
$ws.Range("J131").Value = 934.04346
$ws.Range("L131").Value = 2802.13038
$ws.Range("N131").Value = -12882.13038

# Row 136
$ws.Range("H136").Value = 2808.8262
$ws.Range("I136").Value = 1206.6666
$ws.Range("J136").Value = 3838.7856
$ws.Range("K136").Value = 3619.9998
$ws.Range("L136").Value = 11516.3568
$ws.Range("M136").Value = 1480.0002
$ws.Range("N136").Value = -21716.3568

# Row 137
$ws.Range("H137").Value = 4775.722
$ws.Range("I137").Value = 2039.3125
$ws.Range("J137").Value = 6964.85
$ws.Range("K137").Value = 6117.9375
$ws.Range("L137").Value = 20894.55
$ws.Range("M137").Value = -1017.9375
$ws.Range("N137").Value = -31094.55

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 16692.055
$ws.Range("I113").Value = 2054.4
$ws.Range("J113").Value = 34989.125
$ws.Range("K113").Value = 2054.4
$ws.Range("L113").Value = 34989.125
$ws.Range("M113").Value = 115.5999999999999
$ws.Range("N113").Value = -39329.125

# Row 122
$ws.Range("H122").Value = 27783622
$ws.Range("I122").Value = 7013.75
$ws.Range("J122").Value = 83336830
$ws.Range("K122").Value = 21041.25
$ws.Range("L122").Value = 250010490
$ws.Range("M122").Value = -18591.25
$ws.Range("N122").Value = -250015390

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 38464730
$ws.Range("I22").Value = 480
$ws.Range("J22").Value = 55559948
$ws.Range("K22").Value = 480
$ws.Range("L22").Value = 55559948
$ws.Range("M22").Value = -185
$ws.Range("N22").Value = -55560538

# Row 27
$ws.Range("H27").Value = 38464730
$ws.Range("I27").Value = 480
$ws.Range("J27").Value = 55559948
$ws.Range("K27").Value = 480
$ws.Range("L27").Value = 55559948
$ws.Range("M27").Value = -373
$ws.Range("N27").Value = -55560162

# Row 46
$ws.Range("H46").Value = 62500704
$ws.Range("I46").Value = 841.8
$ws.Range("J46").Value = 166667140
$ws.Range("K46").Value = 841.8
$ws.Range("L46").Value = 166667140
$ws.Range("M46").Value = -653.8
$ws.Range("N46").Value = -166667516

# Row 55
$ws.Range("H55").Value = 12500089
$ws.Range("I55").Value = 20833374
$ws.Range("J55").Value = 162.5
$ws.Range("K55").Value = 20833374
$ws.Range("L55").Value = 162.5
$ws.Range("M55").Value = -20833201
$ws.Range("N55").Value = -508.5

# Row 122
$ws.Range("H122").Value = 9222566
$ws.Range("I122").Value = 1332537.2
$ws.Range("J122").Value = 25002624
$ws.Range("K122").Value = 3997611.6
$ws.Range("L122").Value = 75007872
$ws.Range("M122").Value = -3995161.6
$ws.Range("N122").Value = -75012772

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1759.9714
$ws.Range("I122").Value = 1557.1428
$ws.Range("K122").Value = 4671.428400000001
$ws.Range("M122").Value = -2221.428400000001
